$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 09:03"

# Row 4
$ws.Range("B4").Value = 2462708
$ws.Range("C4").Value = 154
$ws.Range("D4").Value = 1040608
$ws.Range("E4").Value = 1297818
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 124282

# Row 7
$ws.Range("B7").Value = 473719
$ws.Range("C7").Value = 734
$ws.Range("D7").Value = 271723
$ws.Range("E7").Value = 187089
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 14907

# Row 38
$ws.Range("B38").Value = 40008
$ws.Range("C38").Value = 994
$ws.Range("D38").Value = 17758
$ws.Range("E38").Value = 21183
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 16
$ws.Range("H38").Value = 1067

# Row 44
$ws.Range("B44").Value = 30175
$ws.Range("C44").Value = 535
$ws.Range("D44").Value = 10174
$ws.Range("E44").Value = 19326
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 36
$ws.Range("H44").Value = 675

# Row 51
$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 22488
$ws.Range("C51").Value = 771
$ws.Range("D51").Value = 11335
$ws.Range("E51").Value = 10756
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 11
$ws.Range("H51").Value = 397

# Row 52
$ws.Range("A52").Value = "Israel"
$ws.Range("B52").Value = 22044
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 15940
$ws.Range("E52").Value = 5796
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 308

# Row 53
$ws.Range("A53").Value = "Nigeria"
$ws.Range("B53").Value = 22020
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 7613
$ws.Range("E53").Value = 13865
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 542

# Row 84
$ws.Range("B84").Value = 5150
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 2950
$ws.Range("E84").Value = 2081
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 119

# Row 91
$ws.Range("B91").Value = 4123
$ws.Range("C91").Value = 9
$ws.Range("D91").Value = 2640
$ws.Range("E91").Value = 906
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 577

# Row 135
$ws.Range("B135").Value = 917
$ws.Range("C135").Value = 3
$ws.Range("D135").Value = 776
$ws.Range("E135").Value = 127
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 14

# Row 163
$ws.Range("A163").Value = "Siria"
$ws.Range("B163").Value = 242
$ws.Range("C163").Value = 11
$ws.Range("D163").Value = 96
$ws.Range("E163").Value = 139
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 7

# Row 164
$ws.Range("A164").Value = "Martinica"
$ws.Range("B164").Value = 236
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 98
$ws.Range("E164").Value = 124
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 14

# Row 208
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("B208").Value = 13
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 13
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# Row 209
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 211
$ws.Range("A211").Value = "Montserrat"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1

# Row 212
$ws.Range("A212").Value = "Seychelles"
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 11
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0
